$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new value looks like a plain number (e.g. "239.14").
# Excel would auto-convert these to numeric cells, but the source data is
# formatted text (e.g. "1.841.65" style price strings), so force Text format
# before assigning, then clear the temporary formatting back off afterwards.
$numericLookingDCells = @("D4","D5","D6","D7","D8","D9","D10","D11","D13","D14","D15","D16","D17","D19","D20","D22","D24","D25","D26","D28","D29","D30","D31","D32","D33","D35","D37","D38","D40","D41","D42","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($addr in $numericLookingDCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '29.376.53'
$ws.Range("E2").Value = '  -0.41%  '
$ws.Range("D3").Value = '1.840.77'
$ws.Range("E3").Value = '  -0.43%  '
$ws.Range("D4").Value = '0.9999'
$ws.Range("E4").Value = '  +0.16%  '
$ws.Range("D5").Value = '239.14'
$ws.Range("E5").Value = '  -0.54%  '
$ws.Range("D6").Value = '0.6267'
$ws.Range("E6").Value = '  -0.51%  '
$ws.Range("D7").Value = '1.002'
$ws.Range("E7").Value = '  +0.21%  '
$ws.Range("D8").Value = '0.07391'
$ws.Range("E8").Value = '  -1.23%  '
$ws.Range("D9").Value = '0.2889'
$ws.Range("E9").Value = '  -0.80%  '
$ws.Range("D10").Value = '24.82'
$ws.Range("E10").Value = '  +0.75%  '
$ws.Range("D11").Value = '0.07721'
$ws.Range("E11").Value = '  -0.27%  '
$ws.Range("D12").Value = '1.834.11'
$ws.Range("E12").Value = '  -0.71%  '
$ws.Range("D13").Value = '4.959'
$ws.Range("E13").Value = '  -1.14%  '
$ws.Range("D14").Value = '0.6678'
$ws.Range("E14").Value = '  -2.00%  '
$ws.Range("D15").Value = '0.00001036'
$ws.Range("E15").Value = '  -0.81%  '
$ws.Range("D16").Value = '81.47'
$ws.Range("E16").Value = '  -0.93%  '
$ws.Range("D17").Value = '6.245'
$ws.Range("E17").Value = '  -0.15%  '
$ws.Range("D18").Value = '29.394.18'
$ws.Range("E18").Value = '  -0.30%  '
$ws.Range("D19").Value = '234.14'
$ws.Range("E19").Value = '  +2.01%  '
$ws.Range("D20").Value = '12.28'
$ws.Range("E20").Value = '  -0.90%  '
$ws.Range("E21").Value = '  +0.18%  '
$ws.Range("D22").Value = '7.283'
$ws.Range("E22").Value = '  -3.38%  '
$ws.Range("E23").Value = '  +0.16%  '
$ws.Range("D24").Value = '157.05'
$ws.Range("E24").Value = '  -1.49%  '
$ws.Range("D25").Value = '8.459'
$ws.Range("E25").Value = '  -0.70%  '
$ws.Range("D26").Value = '0.1342'
$ws.Range("E26").Value = '  -1.98%  '
$ws.Range("E27").Value = '  -1.40%  '
$ws.Range("D28").Value = '0.07249'
$ws.Range("E28").Value = '  +11.23%  '
$ws.Range("D29").Value = '1.493'
$ws.Range("E29").Value = '  +4.84%  '
$ws.Range("D30").Value = '1.483'
$ws.Range("E30").Value = '  -0.31%  '
$ws.Range("D31").Value = '4.024'
$ws.Range("E31").Value = '  -1.97%  '
$ws.Range("D32").Value = '4.026'
$ws.Range("E32").Value = '  -1.90%  '
$ws.Range("D33").Value = '1.158'
$ws.Range("E33").Value = '  +1.11%  '
$ws.Range("E34").Value = '  -1.04%  '
$ws.Range("D35").Value = '0.7142'
$ws.Range("E35").Value = '  +2.00%  '
$ws.Range("E36").Value = '  +0.22%  '
$ws.Range("D37").Value = '0.01833'
$ws.Range("E37").Value = '  -1.69%  '
$ws.Range("D38").Value = '2.792'
$ws.Range("E38").Value = '  -1.75%  '
$ws.Range("D39").Value = '1.231.63'
$ws.Range("E39").Value = '  -2.70%  '
$ws.Range("D40").Value = '6.776'
$ws.Range("E40").Value = '  -0.84%  '
$ws.Range("D41").Value = '0.9500'
$ws.Range("E41").Value = '  +1.35%  '
$ws.Range("D42").Value = '1.002'
$ws.Range("E42").Value = '  +0.18%  '
$ws.Range("D43").Value = '2.000.31'
$ws.Range("E43").Value = '  -0.91%  '
$ws.Range("D44").Value = '101.17'
$ws.Range("E44").Value = '  -0.16%  '
$ws.Range("D45").Value = '65.18'
$ws.Range("E45").Value = '  -1.60%  '
$ws.Range("D46").Value = '0.00000000116'
$ws.Range("E46").Value = '  -0.77%  '
$ws.Range("D47").Value = '1.694'
$ws.Range("E47").Value = '  -2.77%  '
$ws.Range("D48").Value = '6.949'
$ws.Range("D49").Value = '8.912'
$ws.Range("E49").Value = '  -1.17%  '
$ws.Range("D50").Value = '0.1128'
$ws.Range("E50").Value = '  -3.34%  '
$ws.Range("D51").Value = '0.3877'
$ws.Range("E51").Value = '  -2.18%  '

foreach ($addr in $numericLookingDCells) {
    $ws.Range($addr).ClearFormats()
}
